# Add new cocoa/chocolate-industry processor companies to the supply-chain sheet,
# correct Cocoasource SA's city, format the Volume column with thousands separators,
# and auto-fit the Contact Email / Volume columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row: Cocoasource SA is headquartered in Vevey, not Nyon ---
$ws.Range("D10").Value = "Vevey"

# --- Append new processor companies (rows 11-19) ---
$newRows = @(
    @("Mars Wrigley Confectionery", "Processor", "USA",         "Hackettstown", "contact@mars.com",                  1800000, 40.8529, -74.8299, "Top global confectionery producer"),
    @("Ferrero Group",              "Processor", "Italy",       "Alba",         "info@ferrero.com",                  1300000, 44.6974,   8.0341, "Nutella and premium chocolate"),
    @("Mondelēz International",    "Processor", "USA",         "Chicago",      "investor@mdlz.com",                 1200000, 41.8781, -87.6298, "Owns Cadbury, Milka, Toblerone"),
    @("Meiji Co., Ltd.",            "Processor", "Japan",       "Tokyo",        "info@meiji.com",                    1007500, 35.6895, 139.6917, "Japanese confectionery giant"),
    @("Hershey Co",                 "Processor", "USA",         "Hershey",      "consumerrelations@hersheys.com",     806600, 40.2859, -76.6502, "Major chocolate manufacturer"),
    @("Nestlé S.A.",                "Processor", "Switzerland", "Vevey",        "mediarelations@nestle.com",          763600, 46.4628,   6.8431, "Global food and chocolate producer"),
    @("Lindt & Sprüngli AG",        "Processor", "Switzerland", "Kilchberg",    "contact@lindt.com",                   457400, 47.3277,   8.5517, "Premium Swiss chocolate brand"),
    @("Pladis",                     "Processor", "UK",          "London",       "info@pladisglobal.com",               465500, 51.5072,  -0.1276, "Owns McVitie’s and Godiva (license)"),
    @("Ezaki Glico Co., Ltd.",      "Processor", "Japan",       "Osaka",        "support@glico.com",                   315000, 34.6937, 135.5023, "Pocky and Pretz maker")
)

$r = 11
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r++
}

# --- Apply thousands-separator ("Comma") number format to the Volume column ---
$ws.Range("F1:F19").Style = "Comma"

# --- Auto-fit columns to their (now wider/longer) content ---
$ws.Columns.Item(5).EntireColumn.AutoFit()
$ws.Columns.Item(6).EntireColumn.AutoFit()

$ws.Range("A1").Select()
